$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 ("Liste nach Kategorie filtern" / Pulldown) is now done: mark as
# resolved ("X"), drop the "Problem" note, and hide the row like the other
# already-resolved entries.
$ws.Range("D3").Value = "X"
$ws.Range("E3").ClearContents()
$ws.Rows.Item(3).Hidden = $true

# The autofilter on column D ("Erledigt?") no longer needs to special-case
# the now-gone "x?" value - reapply it with only the blank bucket checked.
$ws.Range("A1:F17").AutoFilter(4, @(""), 7)

# Move the active selection off the edited row.
$ws.Range("A18").Select()
